# Insert a new row at position 377 (shifts existing rows 377:444 down to 378:445)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("377:377").Insert()

# Fill in the values for the newly inserted row 377
$ws.Range("A377").Value = 3
$ws.Range("B377").Value = "Femacal de La Calera"
$ws.Range("C377").Value = "Coquimbo"
$ws.Range("D377").Value2 = 44694
$ws.Range("E377").Value = 5
$ws.Range("F377").Value = 100112003
$ws.Range("G377").Value = "Ajo"
$ws.Range("H377").Value = "Chino"
$ws.Range("I377").Value = "Primera"
$ws.Range("J377").Value = 76
$ws.Range("K377").Value = 19000
$ws.Range("L377").Value = 19500
$ws.Range("M377").Value = 19250
$ws.Range("N377").Value = "$/caja 10 kilos"
$ws.Range("O377").Value = "China"
$ws.Range("P377").Value = 1925
$ws.Range("Q377").Value = 10
$ws.Range("R377").Value = "Hortaliza"
